$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 293; this shifts the former rows 293:446
# down to 294:447, matching the target diff exactly (every row's full
# contents move down by one).
$ws.Rows("293:293").Insert()

# Populate the newly-inserted row 293 with its data (same fixed columns
# used by every data row, plus the new record's own date/price/quality
# figures).
$ws.Range("A293").Value2 = 9
$ws.Range("B293").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C293").Value2 = "Metropolitana"
$ws.Range("D293").Value2 = 45097
$ws.Range("E293").Value2 = 13
$ws.Range("F293").Value2 = 300000001
$ws.Range("G293").Value2 = "Rabanito"
$ws.Range("H293").Value2 = "Sin especificar"
$ws.Range("I293").Value2 = "Primera"
$ws.Range("J293").Value2 = 7000
$ws.Range("K293").Value2 = 3000
$ws.Range("L293").Value2 = 3000
$ws.Range("M293").Value2 = 3000
$ws.Range("N293").Value2 = '$/cien unidades (volumen en unidades)'
$ws.Range("O293").Value2 = "Provincia de Chacabuco"
$ws.Range("P293").Value2 = 30
$ws.Range("Q293").Value2 = 100
$ws.Range("R293").Value2 = "Hortaliza"
